# Rencana penyusunan bukti kelengkapan
# Adds a new "Bukti Kelengkapan" (Asesi) row to each of the relevant
# reference sheets, mirroring the existing NAMA ID / NAMA MENU / JENIS AKUN
# table pattern already present on every sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "{View} ID FORM" -> new row 10 : id_form_bukti
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("{View} ID FORM")
$ws.Range("B10").Value = "id_form_bukti"
$ws.Range("C10").Value = "Bukti Kelengkapan"
$ws.Range("D10").Value = "Asesi"
# Column C on this sheet carries an explicit (non-default) cell style on
# every existing data row -- copy it down from the row above so the new
# row matches.
$ws.Range("C9").Copy()
$ws.Range("C10").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("D11").Select()

# ---------------------------------------------------------------------
# Sheet "{View} ID TABEL" -> new row 12 : id_tabel_bukti
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("{View} ID TABEL")
$ws.Range("B12").Value = "id_tabel_bukti"
$ws.Range("C12").Value = "Bukti Kelengkapan"
$ws.Range("D12").Value = "Asesi"
$ws.Range("B13").Select()

# ---------------------------------------------------------------------
# Sheet "{View} ID MODAL" -> new row 11 : id_modal_form_bukti
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("{View} ID MODAL")
$ws.Range("B11").Value = "id_modal_form_bukti"
$ws.Range("C11").Value = "Bukti Kelengkapan"
$ws.Range("D11").Value = "Asesi"
$ws.Range("B12").Select()

# ---------------------------------------------------------------------
# Sheet "{View} ID PAGE" -> no new row, only selection moved
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("{View} ID PAGE")
$ws.Range("F17").Select()

# ---------------------------------------------------------------------
# Sheet "{Model} DATATABLES" -> new row 12 : datatabel_bukti
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("{Model} DATATABLES")
$ws.Range("B12").Value = "datatabel_bukti"
$ws.Range("C12").Value = "Bukti Kelengkapan"
$ws.Range("D12").Value = "Asesi"
$ws.Range("C15").Select()

# ---------------------------------------------------------------------
# Sheet "{Model} FORM" -> new row 7 : form_bukti
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("{Model} FORM")
$ws.Range("B7").Value = "form_bukti"
$ws.Range("C7").Value = "Bukti Kelengkapan"
$ws.Range("D7").Value = "Asesi"
$ws.Range("D8").Select()

# ---------------------------------------------------------------------
# Sheet "{Controller} aksiTambahData" -> new row 7 : satuData_bukti
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("{Controller} aksiTambahData")
$ws.Range("B7").Value = "satuData_bukti"
$ws.Range("C7").Value = "Bukti Kelengkapan"
$ws.Range("D7").Value = "Asesi"
$ws.Range("B8").Select()

# ---------------------------------------------------------------------
# Sheet "{Controller} aksiAmbilData" -> new rows 24 & 25 :
#   satuData_bukti / datatabel_bukti
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("{Controller} aksiAmbilData")
$ws.Range("B24").Value = "satuData_bukti"
$ws.Range("C24").Value = "Bukti Kelengkapan"
$ws.Range("D24").Value = "Asesi"
$ws.Range("B25").Value = "datatabel_bukti"
$ws.Range("C25").Value = "Bukti Kelengkapan"
$ws.Range("D25").Value = "Asesi"
$ws.Activate()
$ws.Range("C25").Select()

# ---------------------------------------------------------------------
# Sheet "{Controller} aksiUpdateData" -> new row 8 : satuData_bukti
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("{Controller} aksiUpdateData")
$ws.Range("B8").Value = "satuData_bukti"
$ws.Range("C8").Value = "Bukti Kelengkapan"
$ws.Range("D8").Value = "Asesi"
$ws.Range("B9").Select()

# Re-activate the sheet that was selected/active in the original file so
# the saved workbook opens on the same tab.
$wb.Worksheets.Item("{View} ID FORM").Activate()
